$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-25"

# Update the header label in I1 ("2022 (through 06-24)" -> "2022 (through 06-25)")
$ws.Range("I1").Value = "2022 (through 06-25)"

# Update July total (I7): 114 -> 118
$ws.Range("I7").Value = 118

# Update grand total (I14): 777 -> 781
$ws.Range("I14").Value = 781
